$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the bad-pair annotations that are no longer needed (E18, E19)
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()

# Set the selection to match the saved view state
$ws.Range("E18").Select()

# Autofit column E so its width matches the longest remaining entry
$ws.Columns("E:E").AutoFit() | Out-Null
